$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 133.9646796666667
$ws.Range("H2").Value = 401.894039
$ws.Range("I2").Value = 0.2795129415517746
$ws.Range("J2").Value = 0.2795129415517745
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 3.090355666666666
$ws.Range("N2").Value = 9.271066999999999
$ws.Range("O2").Value = 0.06928583878088775
$ws.Range("P2").Value = 0.06928583878088775
$ws.Range("Q2").Value = 413.9985069410681
$ws.Range("R2").Value = 3725.986562469613
$ws.Range("S2").Value = 0.01936628860552796
$ws.Range("T2").Value = 0.01936628860552795
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 133.9646796666667
$ws.Range("H3").Value = 401.894039
$ws.Range("I3").Value = 0.2795129415517746
$ws.Range("J3").Value = 0.2795129415517745
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 25.17096033333333
$ws.Range("N3").Value = 75.512881
$ws.Range("O3").Value = 0.5643334579338453
$ws.Range("P3").Value = 0.5643334579338454
$ws.Range("Q3").Value = 3372.019637957373
$ws.Range("R3").Value = 30348.17674161636
$ws.Range("S3").Value = 0.1577385048431738
$ws.Range("T3").Value = 0.1577385048431737
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 133.9646796666667
$ws.Range("H4").Value = 401.894039
$ws.Range("I4").Value = 0.2795129415517746
$ws.Range("J4").Value = 0.2795129415517745
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 16.34167533333333
$ws.Range("N4").Value = 49.025026
$ws.Range("O4").Value = 0.366380703285267
$ws.Range("P4").Value = 0.366380703285267
$ws.Range("Q4").Value = 2189.207301246668
$ws.Range("R4").Value = 19702.86571122001
$ws.Range("S4").Value = 0.1024081481030729
$ws.Range("T4").Value = 0.1024081481030729
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 276.4348856666666
$ws.Range("H5").Value = 829.3046569999999
$ws.Range("I5").Value = 0.576772386814763
$ws.Range("J5").Value = 0.576772386814763
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.090355666666666
$ws.Range("N5").Value = 9.271066999999999
$ws.Range("O5").Value = 0.06928583878088775
$ws.Range("P5").Value = 0.06928583878088775
$ws.Range("Q5").Value = 854.2821153843352
$ws.Range("R5").Value = 7688.539038459017
$ws.Range("S5").Value = 0.03996215860611549
$ws.Range("T5").Value = 0.03996215860611549
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 276.4348856666666
$ws.Range("H6").Value = 829.3046569999999
$ws.Range("I6").Value = 0.576772386814763
$ws.Range("J6").Value = 0.576772386814763
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 25.17096033333333
$ws.Range("N6").Value = 75.512881
$ws.Range("O6").Value = 0.5643334579338453
$ws.Range("P6").Value = 0.5643334579338454
$ws.Range("Q6").Value = 6958.131541865199
$ws.Range("R6").Value = 62623.1838767868
$ws.Range("S6").Value = 0.3254919554919325
$ws.Range("T6").Value = 0.3254919554919326
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 276.4348856666666
$ws.Range("H7").Value = 829.3046569999999
$ws.Range("I7").Value = 0.576772386814763
$ws.Range("J7").Value = 0.576772386814763
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.34167533333333
$ws.Range("N7").Value = 49.025026
$ws.Range("O7").Value = 0.366380703285267
$ws.Range("P7").Value = 0.366380703285267
$ws.Range("Q7").Value = 4517.409152371785
$ws.Range("R7").Value = 40656.68237134608
$ws.Range("S7").Value = 0.2113182727167149
$ws.Range("T7").Value = 0.2113182727167149
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 68.87942233333334
$ws.Range("H8").Value = 206.638267
$ws.Range("I8").Value = 0.1437146716334625
$ws.Range("J8").Value = 0.1437146716334625
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 3.090355666666666
$ws.Range("N8").Value = 9.271066999999999
$ws.Range("O8").Value = 0.06928583878088775
$ws.Range("P8").Value = 0.06928583878088775
$ws.Range("Q8").Value = 212.8619131245432
$ws.Range("R8").Value = 1915.757218120889
$ws.Range("S8").Value = 0.009957391569244307
$ws.Range("T8").Value = 0.009957391569244306
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 68.87942233333334
$ws.Range("H9").Value = 206.638267
$ws.Range("I9").Value = 0.1437146716334625
$ws.Range("J9").Value = 0.1437146716334625
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 25.17096033333333
$ws.Range("N9").Value = 75.512881
$ws.Range("O9").Value = 0.5643334579338453
$ws.Range("P9").Value = 0.5643334579338454
$ws.Range("Q9").Value = 1733.761207335247
$ws.Range("R9").Value = 15603.85086601723
$ws.Range("S9").Value = 0.08110299759873903
$ws.Range("T9").Value = 0.08110299759873903
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 68.87942233333334
$ws.Range("H10").Value = 206.638267
$ws.Range("I10").Value = 0.1437146716334625
$ws.Range("J10").Value = 0.1437146716334625
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 16.34167533333333
$ws.Range("N10").Value = 49.025026
$ws.Range("O10").Value = 0.366380703285267
$ws.Range("P10").Value = 0.366380703285267
$ws.Range("Q10").Value = 1125.605156918882
$ws.Range("R10").Value = 10130.44641226994
$ws.Range("S10").Value = 0.05265428246547921
$ws.Range("T10").Value = 0.05265428246547921
